$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H9").Value = 0.1146
$ws.Range("I9").Value = -0.2406
$ws.Range("H13").Value = 0.0006
$ws.Range("I13").Value = 0.0491
$ws.Range("J13").Value = 0.0711
$ws.Range("K13").Value = 0.0484
$ws.Range("L13").Value = 0.1798
$ws.Range("M13").Value = 0.0433
$ws.Range("N13").Value = -0.0098
$ws.Range("O13").Value = -0.0532
$ws.Range("P13").Value = -0.0965
$ws.Range("H17").Value = 0.0682
$ws.Range("I17").Value = -0.2077
$ws.Range("J17").Value = -0.1207
$ws.Range("K17").Value = -0.0848
$ws.Range("L17").Value = -0.0507
$ws.Range("M17").Value = -0.0232
$ws.Range("N17").Value = -0.0329
$ws.Range("O17").Value = -0.0644
$ws.Range("P17").Value = -0.0547
$ws.Range("H23").Value = -0.0884
$ws.Range("I23").Value = -0.7622
$ws.Range("J23").Value = -0.4149
$ws.Range("K23").Value = -0.27
$ws.Range("L23").Value = -0.1292
$ws.Range("M23").Value = -0.0661
$ws.Range("N23").Value = -0.1769
$ws.Range("H25").Value = -0.1803
$ws.Range("I25").Value = -0.1796
$ws.Range("J25").Value = -0.0664
$ws.Range("K25").Value = -0.0543
$ws.Range("L25").Value = -0.601
$ws.Range("M25").Value = -0.633
$ws.Range("N25").Value = -0.3094
$ws.Range("O25").Value = -0.1899
$ws.Range("P25").Value = -0.3258
$ws.Range("G26").Value = -0.9899
$ws.Range("H27").Value = -1.4822
$ws.Range("I27").Value = -1.019
$ws.Range("J27").Value = -0.4185
$ws.Range("K27").Value = -0.4068
$ws.Range("L27").Value = -0.3092
$ws.Range("M27").Value = -0.1444
$ws.Range("N27").Value = -0.0072
$ws.Range("O27").Value = -0.0016
$ws.Range("P27").Value = 0
$ws.Range("H29").Value = -4.1061
$ws.Range("I29").Value = -3.501
$ws.Range("J29").Value = -2.2257
$ws.Range("K29").Value = -2.1498
$ws.Range("L29").Value = -2.9398
$ws.Range("M29").Value = -1.3519
$ws.Range("N29").Value = -0.9488
$ws.Range("O29").Value = -0.4768
$ws.Range("P29").Value = -0.5114
$ws.Range("H33").Value = 0.1018
$ws.Range("I33").Value = -0.1918
$ws.Range("H41").Value = 0.0154
$ws.Range("I41").Value = 0.3123
$ws.Range("H45").Value = -0.0492
$ws.Range("I45").Value = -0.154
$ws.Range("J45").Value = -0.105
$ws.Range("K45").Value = 0.0207
$ws.Range("L45").Value = 0.0393
$ws.Range("M45").Value = 0.1934
$ws.Range("N45").Value = 0.1558
$ws.Range("O45").Value = 0.1266
$ws.Range("P45").Value = 0.0979
$ws.Range("G54").Value = -0.1281
$ws.Range("H55").Value = -0.2303
$ws.Range("I55").Value = -0.0825
$ws.Range("J55").Value = -0.0427
$ws.Range("K55").Value = -0.0331
$ws.Range("L55").Value = -0.0219
$ws.Range("M55").Value = -0.0136
$ws.Range("O55").Value = 0.0013
$ws.Range("P55").Value = 0.0014
$ws.Range("H63").Value = 0.3384
$ws.Range("I63").Value = -0.0499
$ws.Range("H67").Value = -0.0409
$ws.Range("I67").Value = 0.0632
$ws.Range("J67").Value = 0.1552
$ws.Range("K67").Value = 0.3498
$ws.Range("L67").Value = 0.4344
$ws.Range("M67").Value = 0.2602
$ws.Range("N67").Value = 0.1816
$ws.Range("O67").Value = 0.0138
$ws.Range("P67").Value = -0.0959
$ws.Range("H71").Value = 0.0699
$ws.Range("I71").Value = -0.0155
$ws.Range("J71").Value = 0.0034
$ws.Range("K71").Value = -0.0213
$ws.Range("L71").Value = 0.002
$ws.Range("M71").Value = 0.0011
$ws.Range("N71").Value = 0.0016
$ws.Range("O71").Value = 0.0029
$ws.Range("P71").Value = 0.0025
$ws.Range("H77").Value = 0.1253
$ws.Range("I77").Value = 0.3485
$ws.Range("J77").Value = 0.3561
$ws.Range("K77").Value = 0.0952
$ws.Range("L77").Value = 0.0188
$ws.Range("M77").Value = -0.175
$ws.Range("N77").Value = -0.1737
$ws.Range("H79").Value = 0.016
$ws.Range("I79").Value = 0.0025
$ws.Range("J79").Value = 0.0021
$ws.Range("K79").Value = 0.0023
$ws.Range("L79").Value = 0.021
$ws.Range("M79").Value = 0.025
$ws.Range("N79").Value = 0.0136
$ws.Range("O79").Value = 0.0089
$ws.Range("P79").Value = 0.0143
$ws.Range("G80").Value = 0.0022
$ws.Range("H81").Value = -0.0092
$ws.Range("I81").Value = 0.0099
$ws.Range("J81").Value = 0.002
$ws.Range("K81").Value = 0.0043
$ws.Range("L81").Value = 0.0058
$ws.Range("M81").Value = 0.0029
$ws.Range("N81").Value = -0.0023
$ws.Range("O81").Value = -0.0021
$ws.Range("P81").Value = -0.0002
$ws.Range("H83").Value = -0.4026
$ws.Range("I83").Value = -0.2281
$ws.Range("J83").Value = -0.1297
$ws.Range("K83").Value = -0.1827
$ws.Range("L83").Value = 0.1122
$ws.Range("M83").Value = -0.2208
$ws.Range("N83").Value = -0.323
$ws.Range("O83").Value = 1.3531
$ws.Range("P83").Value = 0.6281
$ws.Range("H87").Value = 0.229
$ws.Range("I87").Value = -0.1419
$ws.Range("H95").Value = -0.0851
$ws.Range("I95").Value = 0.2358
$ws.Range("H99").Value = -0.139
$ws.Range("I99").Value = -0.152
$ws.Range("J99").Value = -0.1135
$ws.Range("K99").Value = -0.0122
$ws.Range("L99").Value = 0.0888
$ws.Range("M99").Value = 0.2859
$ws.Range("N99").Value = 0.2279
$ws.Range("O99").Value = 0.1712
$ws.Range("P99").Value = 0.1156
$ws.Range("G108").Value = 0
$ws.Range("H109").Value = -0.0039
$ws.Range("I109").Value = 0.0031
$ws.Range("J109").Value = 0.0023
$ws.Range("K109").Value = 0.0025
$ws.Range("L109").Value = 0.0033
$ws.Range("M109").Value = 0.0025
$ws.Range("O109").Value = 0.0012
$ws.Range("P109").Value = 0.0016
